$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 'Larvik/Norge'
$ws.Range("B6").Value = 'Hirtshals/Danmark'
$ws.Range("C6").Value = 'mandag, onsdag, fredag '
$ws.Range("D6").Value = '08:00, 17:00'
$ws.Range("E6").Value = 'NOK 549.00'
$ws.Range("F6").Value = 'NOK 399.00'
$ws.Range("G6").Value = 'NOK 100.00'
$ws.Range("H6").Value = 'NOK 149.00'

$ws.Range("A7").Value = 'Oslo/Norge'
$ws.Range("B7").Value = 'København/Danmark'
$ws.Range("C7").Value = 'fredag, lørdag, søndag, mandag'
$ws.Range("D7").Value = 0.58333333333333337
$ws.Range("D7").NumberFormat = "h:mm"
$ws.Range("E7").Value = 'NOK 499.00'
$ws.Range("F7").Value = 'NOK 299.00'
$ws.Range("G7").Value = 'NOK 100.00'
$ws.Range("H7").Value = 'NOK 149.00'

$ws.Range("A8").Value = 'Bodø/Norge'
$ws.Range("B8").Value = 'Moskenes/Norge'
$ws.Range("C8").Value = 'mandag, torsdag, lørdag'
$ws.Range("D8").Value = '12:00, 22:00'
$ws.Range("E8").Value = 'NOK 749.00'
$ws.Range("F8").Value = 'NOK 299.00'
$ws.Range("G8").Value = 'NOK 100.00'
$ws.Range("H8").Value = 'NOK 0.00'

$ws.Range("A9").Value = 'Sandefjord/Norge'
$ws.Range("B9").Value = 'Strømstad/Sverige'
$ws.Range("C9").Value = 'mandag, onsdag, fredag, søndag'
$ws.Range("D9").Value = '09:00, 16:00'
$ws.Range("E9").Value = 'NOK 399.00'
$ws.Range("F9").Value = 'NOK 299.00'
$ws.Range("G9").Value = 'NOK 100.00'
$ws.Range("H9").Value = 'NOK 99.00'

$ws.Range("A10").Value = 'Stavanger/Norge'
$ws.Range("B10").Value = 'Bergen/Norge'
$ws.Range("C10").Value = 'mandag, tirsdag, onsdag, torsdag, fredag'
$ws.Range("D10").Value = '07:00, 15:00'
$ws.Range("E10").Value = 'NOK 499.00'
$ws.Range("F10").Value = 'NOK 349.00'
$ws.Range("G10").Value = 'NOK 100.00'
$ws.Range("H10").Value = 'NOK 199.00'

$ws.Range("A11").Value = 'Oslo/Norge'
$ws.Range("B11").Value = 'Kiel/Tyskland'
$ws.Range("C11").Value = 'tirsdag, torsdag, lørdag'
$ws.Range("D11").Value = 0.58333333333333337
$ws.Range("D11").NumberFormat = "h:mm"
$ws.Range("E11").Value = 'NOK 399.00'
$ws.Range("F11").Value = 'NOK 349.00'
$ws.Range("G11").Value = 'NOK 100.00'
$ws.Range("H11").Value = 'NOK 299.00'

$ws.Range("L20").Select()
